$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dateCreated encoding from iso8601 to w3cdtf
$ws.Range("S1").Value = '<mods:originInfo><mods:dateCreated encoding="w3cdtf">'

# Reflect the view state change recorded in the saved file: the user's
# active cell/selection moved to S1 (scrolling the window so later columns,
# like H onward, are in view) while editing that cell.
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("S1").Select()
